$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 52: D52 already held "Indie Project: Cleaning up JSP/Servlet connections..." text (shared string moved
# index but text is unchanged) -- nothing to change there content-wise.

# Row 53: becomes a real data row: date 2019-03-30 (serial 43554), row height 30,
# D53 gets the new task description.
$ws.Range("A53").Value = (Get-Date -Year 2019 -Month 3 -Day 30 -Hour 0 -Minute 0 -Second 0).Date
$ws.Rows.Item(53).RowHeight = 30
$ws.Range("D53").Value = "Indie Project: Fixed IncomeSkew mapping, added new test for surveys to check details explicitly.  Added login servlet so I could reinstate the login nav option."

# Row 54: remove old D54 content ("Uh-oh, survey not mapped right or servlet doesn't display right: ")
$ws.Range("D54").Clear()

# Row 55: remove old D55 content ("Impact of other financial factors: Unmet goals caused frustration.")
$ws.Range("D55").Clear()

# Update the view's selection to match the new state (D54 is now the active cell).
$ws.Range("D54").Select()
